$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 294; existing rows 294-382 shift down to 295-383.
$ws.Rows("294:294").Insert()

# Populate the newly inserted row 294 with the new data point.
$ws.Range("A294").Value = 11
$ws.Range("B294").Value = "Vega Monumental Concepción"
$ws.Range("C294").Value = "Bíobío"
$ws.Range("D294").Value = 44642
$ws.Range("E294").Value = 8
$ws.Range("F294").Value = 100112020
$ws.Range("G294").Value = "Tomate"
$ws.Range("H294").Value = "Larga vida"
$ws.Range("I294").Value = "Primera"
$ws.Range("J294").Value = 220
$ws.Range("K294").Value = 14000
$ws.Range("L294").Value = 15000
$ws.Range("M294").Value = 14545
$ws.Range("N294").Value = "`$/bandeja 18 kilos"
$ws.Range("O294").Value = "Provincia de Quillota"
$ws.Range("P294").Value = 808
$ws.Range("Q294").Value = 18
$ws.Range("R294").Value = "Hortaliza"
